# Insert a new row above row 708, duplicating the adjacent row's format/content
# (so the new date cell keeps its original text representation instead of being
# re-interpreted as a real date), then overwrite the few cells that actually
# differ for the new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A708").EntireRow.Insert()
$ws.Range("A707:D707").Copy($ws.Range("A708:D708"))

$ws.Range("C708").Value = 5
